$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1764
$ws.Range("F3").Value = 10318
$ws.Range("F6").Value = 603
$ws.Range("F8").Value = 1675
$ws.Range("F9").Value = 416
$ws.Range("F11").Value = 231
$ws.Range("F13").Value = 508
$ws.Range("F16").Value = 34
$ws.Range("F18").Value = 31
$ws.Range("F19").Value = 110
$ws.Range("F20").Value = 378
$ws.Range("F23").Value = 27
$ws.Range("F25").Value = 1179
$ws.Range("C26").Value = '上海·cdc动漫展'
$ws.Range("E26").Value = '2024.06.29 10:00-06.30 17:00'
$ws.Range("F26").Value = 203
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=85110'
$ws.Range("I26").Value = '//i2.hdslb.com/bfs/openplatform/202405/RMpaP6sF1714725969882.jpeg'
$ws.Range("C27").Value = '上海·创世次元动漫游戏嘉年华3.0'
$ws.Range("D27").Value = '中环立交桥苏宁天御国际广场西南侧约240米 轮客行轮滑馆(普陀店)'
$ws.Range("F27").Value = 397
$ws.Range("G27").Value = 58
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=86506'
$ws.Range("I27").Value = '//i2.hdslb.com/bfs/openplatform/202405/Clkfdwic1716894666596.jpeg'
$ws.Range("C28").Value = '上海·葬送的芙莉莲ONLY'
$ws.Range("D28").Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Range("E28").Value = '2024.06.29 10:00-06.29 17:00'
$ws.Range("F28").Value = 253
$ws.Range("G28").Value = 65
$ws.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=85193'
$ws.Range("I28").Value = '//i2.hdslb.com/bfs/openplatform/202404/VIM2lfxY1714361685906.jpeg'
$ws.Range("B29").Value = '2024-06-30'
$ws.Range("C29").Value = '上海 星芒旋转 anikura动漫嘉年华'
$ws.Range("D29").Value = '海潮路133号B1 JUMP工坊'
$ws.Range("E29").Value = '2024.06.30 17:00-06.30 21:00'
$ws.Range("F29").Value = 27
$ws.Range("G29").Value = 60
$ws.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=86357'
$ws.Range("I29").Value = '//i2.hdslb.com/bfs/openplatform/202405/RAjlTt6f1716788945862.jpeg'
$ws.Range("F31").Value = 583
$ws.Range("F32").Value = 241
$ws.Range("F33").Value = 738
$ws.Range("F35").Value = 730
$ws.Range("F36").Value = 106162
$ws.Range("F37").Value = 783
$ws.Range("F39").Value = 1271
$ws.Range("F40").Value = 824
$ws.Range("F41").Value = 717
$ws.Range("F42").Value = 749
$ws.Range("F44").Value = 34
$ws.Range("F45").Value = 719
$ws.Range("F47").Value = 713

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 91
$ws.Range("F18").Value = 1124
$ws.Range("F20").Value = 1828
$ws.Range("F21").Value = 1828
$ws.Range("F23").Value = 339
$ws.Range("F25").Value = 88
$ws.Range("F29").Value = 373

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 823
$ws.Range("F6").Value = 2553
$ws.Range("F7").Value = 4199
$ws.Range("F8").Value = 72
$ws.Range("F10").Value = 380
$ws.Range("F11").Value = 252
$ws.Range("F12").Value = 235

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1764
$ws.Range("F4").Value = 823
$ws.Range("F7").Value = 10318
$ws.Range("F9").Value = 4199
$ws.Range("F10").Value = 72
$ws.Range("F11").Value = 380
$ws.Range("F12").Value = 603
$ws.Range("F13").Value = 252
$ws.Range("F14").Value = 1675
$ws.Range("F15").Value = 416
$ws.Range("F16").Value = 231
$ws.Range("F21").Value = 31
$ws.Range("F22").Value = 110
$ws.Range("F23").Value = 1124
$ws.Range("F24").Value = 378
$ws.Range("F26").Value = 27
$ws.Range("F27").Value = 1828
$ws.Range("F29").Value = 1179
$ws.Range("F30").Value = 88
$ws.Range("F32").Value = 253
$ws.Range("F33").Value = 373
$ws.Range("F34").Value = 583
$ws.Range("F36").Value = 738
$ws.Range("F38").Value = 730
$ws.Range("F41").Value = 783
$ws.Range("F43").Value = 824
$ws.Range("F44").Value = 717
$ws.Range("F45").Value = 749
$ws.Range("F48").Value = 719
$ws.Range("F49").Value = 713
